$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 5; $r -le 24; $r++) {
    $ws.Rows.Item($r).RowHeight = 14.25
}
